# Weekly fruit/veggie price update:
# Insert a new data row at row 68 (pushing the existing rows 68-77 down to
# 69-78) and populate the new row with this week's reading for
# "Poroto verde" at Vega Modelo de Temuco.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 68:77 down to 69:78, leaving a blank row 68 for the new entry.
$ws.Rows.Item(68).Insert()

$ws.Range("A68").Value = 10
$ws.Range("B68").Value = "Vega Modelo de Temuco"
$ws.Range("C68").Value = "La Araucanía"
$ws.Range("D68").Value = 44476
$ws.Range("E68").Value = 9
$ws.Range("F68").Value = 100112031
$ws.Range("G68").Value = "Poroto verde"
$ws.Range("H68").Value = "Sin especificar"
$ws.Range("I68").Value = "Primera"
$ws.Range("J68").Value = 20
$ws.Range("K68").Value = 40000
$ws.Range("L68").Value = 40000
$ws.Range("M68").Value = 40000
$ws.Range("N68").Value = "$/malla 25 kilos"
$ws.Range("O68").Value = "Provincia de Limarí"
$ws.Range("P68").Value = 1600
$ws.Range("Q68").Value = 25
$ws.Range("R68").Value = "Hortaliza"
